$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "42.143.20"
Set-TextValue $ws.Range("E2") "  -3.89%  "

Set-TextValue $ws.Range("D3") "2.239.54"
Set-TextValue $ws.Range("E3") "  -4.61%  "

Set-TextValue $ws.Range("E4") "  -0.37%  "

Set-TextValue $ws.Range("D5") "244.78"
Set-TextValue $ws.Range("E5") "  +2.13%  "

Set-TextValue $ws.Range("D6") "0.630"
Set-TextValue $ws.Range("E6") "  -5.38%  "

Set-TextValue $ws.Range("D7") "68.72"
Set-TextValue $ws.Range("E7") "  -5.67%  "

Set-TextValue $ws.Range("E8") "  -0.02%  "

Set-TextValue $ws.Range("D9") "0.552"
Set-TextValue $ws.Range("E9") "  -6.76%  "

Set-TextValue $ws.Range("D10") "0.0983"
Set-TextValue $ws.Range("E10") "  -2.30%  "

Set-TextValue $ws.Range("D11") "58.87"
Set-TextValue $ws.Range("E11") "  -0.76%  "

Set-TextValue $ws.Range("D12") "36.02"
Set-TextValue $ws.Range("E12") "  +9.78%  "

Set-TextValue $ws.Range("E13") "  -2.66%  "

Set-TextValue $ws.Range("D14") "6.70"
Set-TextValue $ws.Range("E14") "  -7.53%  "

Set-TextValue $ws.Range("D15") "2.572.04"
Set-TextValue $ws.Range("E15") "  -4.69%  "

Set-TextValue $ws.Range("D16") "14.93"
Set-TextValue $ws.Range("E16") "  -7.20%  "

Set-TextValue $ws.Range("D17") "0.860"
Set-TextValue $ws.Range("E17") "  -4.79%  "

Set-TextValue $ws.Range("D18") "2.244.74"
Set-TextValue $ws.Range("E18") "  -4.52%  "

Set-TextValue $ws.Range("D19") "42.037.36"
Set-TextValue $ws.Range("E19") "  -4.10%  "

Set-TextValue $ws.Range("D20") "0.0₃0964"
Set-TextValue $ws.Range("E20") "  -5.54%  "

Set-TextValue $ws.Range("D21") "73.06"
Set-TextValue $ws.Range("E21") "  -7.04%  "

Set-TextValue $ws.Range("D22") "6.21"
Set-TextValue $ws.Range("E22") "  -7.06%  "

Set-TextValue $ws.Range("D23") "235.41"
Set-TextValue $ws.Range("E23") "  -6.43%  "

Set-TextValue $ws.Range("D24") "2.05"
Set-TextValue $ws.Range("E24") "  +11.75%  "

Set-TextValue $ws.Range("E25") "  -0.07%  "

Set-TextValue $ws.Range("E26") "  -3.49%  "

Set-TextValue $ws.Range("D27") "2.46"
Set-TextValue $ws.Range("E27") "  -0.81%  "

Set-TextValue $ws.Range("D28") "2.24"
Set-TextValue $ws.Range("E28") "  -3.57%  "

Set-TextValue $ws.Range("D29") "9.94"
Set-TextValue $ws.Range("E29") "  -4.42%  "

Set-TextValue $ws.Range("D30") "172.07"
Set-TextValue $ws.Range("E30") "  -2.46%  "

Set-TextValue $ws.Range("D31") "20.49"
Set-TextValue $ws.Range("E31") "  -7.78%  "

Set-TextValue $ws.Range("E32") "  -4.33%  "

Set-TextValue $ws.Range("D33") "0.126"
Set-TextValue $ws.Range("E33") "  -5.32%  "

Set-TextValue $ws.Range("D34") "0.0714"
Set-TextValue $ws.Range("E34") "  -3.90%  "

Set-TextValue $ws.Range("D35") "5.25"
Set-TextValue $ws.Range("E35") "  -1.64%  "

Set-TextValue $ws.Range("D36") "4.69"
Set-TextValue $ws.Range("E36") "  -7.53%  "

Set-TextValue $ws.Range("D37") "3.78"
Set-TextValue $ws.Range("E37") "  +0.57%  "

Set-TextValue $ws.Range("D38") "22.76"
Set-TextValue $ws.Range("E38") "  +21.15%  "

Set-TextValue $ws.Range("E39") "  +5.07%  "

Set-TextValue $ws.Range("D40") "2.30"
Set-TextValue $ws.Range("E40") "  -3.18%  "

Set-TextValue $ws.Range("D41") "5.85"
Set-TextValue $ws.Range("E41") "  -8.48%  "

Set-TextValue $ws.Range("D42") "66.62"
Set-TextValue $ws.Range("E42") "  +2.42%  "

Set-TextValue $ws.Range("D43") "9.16"
Set-TextValue $ws.Range("E43") "  -0.78%  "

Set-TextValue $ws.Range("E44") "  -13.05%  "

Set-TextValue $ws.Range("E45") "  -3.14%  "

Set-TextValue $ws.Range("B46") "Algorand"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D46") "0.190"
Set-TextValue $ws.Range("E46") "  -2.38%  "

Set-TextValue $ws.Range("B47") "BinanceUSD"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D47") "1.00"
Set-TextValue $ws.Range("E47") "  +0.15%  "

Set-TextValue $ws.Range("D48") "4.53"
Set-TextValue $ws.Range("E48") "  +7.59%  "

Set-TextValue $ws.Range("D49") "1.19"
Set-TextValue $ws.Range("E49") "  -2.28%  "

Set-TextValue $ws.Range("B50") "Celestia"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D50") "9.99"
Set-TextValue $ws.Range("E50") "  +6.94%  "

Set-TextValue $ws.Range("D51") "2.31"
Set-TextValue $ws.Range("E51") "  -3.68%  "
